# Edit script: apply DataTier.Net 8 Planning changes (v4.4.1 - removes binding options)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate header (E1) and all data rows/columns ---
$ws.Range("A1").Value = "Title"
$ws.Range("B1").Value = "Details"
$ws.Range("C1").Value = "Difficulty"
$ws.Range("D1").Value = "Status"
$ws.Range("E1").Value = "Notes"
$ws.Range("A2").Value = "Remove Data Watcher"
$ws.Range("B2").Value = "When I first learned Blazor, I needed this. I don't think is needed now."
$ws.Range("C2").Value = "Medium - Took out the DataGateway code"
$ws.Range("D2").Value = "Done Except Docs & Testing"
$ws.Range("E2").Value = "Make sure to update the User's Guide"
$ws.Range("A3").Value = "Remove Binding"
$ws.Range("B3").Value = "Take out code for AllowBinding and Code Generation of ItemCallback"
$ws.Range("C3").Value = "Medium - Effects code generation and GUI"
$ws.Range("D3").Value = "Done Except Docs & Testing"
$ws.Range("E3").Value = "Required update to DataJuggler.Net (.NET Framework version)"
$ws.Range("A4").Value = "Replace Stored Procedure Manager"
$ws.Range("B4").Value = "Stored Procs were too object oriented. Replace with XML loaded."
$ws.Range("C4").Value = "Major - Remove code generated stored procs and builders."
$ws.Range("D4").Value = "Planning"
$ws.Range("E4").Value = "Template Change Also"
$ws.Range("A5").Value = "Remove App Logic Component"
$ws.Range("B5").Value = "Move Connection to DAC"
$ws.Range("C5").Value = "Medium - Move everything to DAC - Investigate removing Data Operations"
$ws.Range("D5").Value = "Planning"
$ws.Range("E5").Value = "Template Change Also"
$ws.Range("A6").Value = "Move Gateway to DAC"
$ws.Range("B6").Value = "This gets it down to 2 projects"
$ws.Range("C6").Value = "Medium - Code Generation and Custom Methods has to change"
$ws.Range("D6").Value = "Planning"
$ws.Range("E6").Value = "Template Change Also"
$ws.Range("A7").Value = "Move Object Library To DAC"
$ws.Range("B7").Value = "This gets it down 1 project"
$ws.Range("C7").Value = "Medium"
$ws.Range("D7").Value = "Planning"
$ws.Range("E7").Value = "Template Change Also"
$ws.Range("A8").Value = "Visual Studio Project Updater"
$ws.Range("B8").Value = "Redesign form and control for 1 project"
$ws.Range("C8").Value = "Medium"
$ws.Range("D8").Value = "Planning"
$ws.Range("A9").Value = "Update Documentation"
$ws.Range("B9").Value = "Lots of the Users Guide and Quick Start Will Change"
$ws.Range("C9").Value = "Medium"
$ws.Range("D9").Value = "Planning"
$ws.Range("A10").Value = "Benchmarks"
$ws.Range("B10").Value = "Test if there are any performance gaines"
$ws.Range("C10").Value = "Easy"
$ws.Range("D10").Value = "Planning"
$ws.Range("A11").Value = "Update Project Template"
$ws.Range("B11").Value = "This goes from a multi-project template to single"
$ws.Range("C11").Value = "Medium"
$ws.Range("D11").Value = "Planning"
$ws.Range("E11").Value = "Build Copy should handle this"
$ws.Range("A12").Value = "Change Custom Method Editors"
$ws.Range("B12").Value = "Use XML and make sure build rewrites base and custom methods"
$ws.Range("C12").Value = "Major "
$ws.Range("D12").Value = "Planning"
$ws.Range("A13").Value = "New Install Release"
$ws.Range("B13").Value = "Update the installer"
$ws.Range("C13").Value = "Easy"
$ws.Range("D13").Value = "Planning"
$ws.Range("A14").Value = "Fix Indent for Custom Methods"
$ws.Range("B14").Value = "Indent is 1 tab to indented."
$ws.Range("C14").Value = "Easy"
$ws.Range("D14").Value = "Done Except Testing"
$ws.Range("E14").Value = "VS 2022 indent doesn't like indented regions"
$ws.Range("A15").Value = "Remove Order By Descending"
$ws.Range("B15").Value = "Field Set "

# --- Make E1 header bold to match existing header style ---
$ws.Range("E1").Font.Bold = $true

# --- Apply highlight fill (theme color 3) to the two "done" rows (A2:A3) ---
$ws.Range("A2:A3").Interior.ThemeColor = 3

# --- Apply highlight fill to row 14 (A14:C14 and E14, but not D14) ---
$ws.Range("A14:C14").Interior.ThemeColor = 3
$ws.Range("E14").Interior.ThemeColor = 3

# --- Column widths (auto-fit approximation) ---
$ws.Columns.Item(1).ColumnWidth = 28.125
$ws.Columns.Item(3).ColumnWidth = 71.125
$ws.Columns.Item(4).ColumnWidth = 21.875
$ws.Columns.Item(5).ColumnWidth = 49.0

# --- View state: scroll so row 4 is at top, select B20 ---
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("B20").Select()
